# Updated cryptos list values (Price / Volume(1h)) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.974.86'
$ws.Range("E2").Value = '  -3.85%  '
$ws.Range("D3").Value = '3.352.50'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.05'
$ws.Range("E5").Value = '  -3.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.51'
$ws.Range("E6").Value = '  -4.92%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -1.74%  '
$ws.Range("E9").Value = '  -3.64%  '
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.403'
$ws.Range("E11").Value = '  -4.58%  '
$ws.Range("D12").Value = '3.930.54'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.23'
$ws.Range("E14").Value = '  -5.39%  '
$ws.Range("D15").Value = '66.967.77'
$ws.Range("E15").Value = '  -3.85%  '
$ws.Range("E16").Value = '  -2.61%  '
$ws.Range("D17").Value = '3.341.12'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '437.75'
$ws.Range("E18").Value = '  -2.97%  '
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.69'
$ws.Range("E20").Value = '  -2.78%  '
$ws.Range("E21").Value = '  -2.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.90'
$ws.Range("E22").Value = '  -1.41%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.522'
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("E25").Value = '  -3.05%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.13'
$ws.Range("E27").Value = '  -4.74%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.98'
$ws.Range("E30").Value = '  -1.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.34'
$ws.Range("E31").Value = '  -5.66%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("E33").Value = '  -3.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.82'
$ws.Range("E34").Value = '  -3.11%  '
$ws.Range("E35").Value = '  -2.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.55'
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.88'
$ws.Range("E37").Value = '  +2.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.86'
$ws.Range("E38").Value = '  -5.24%  '
$ws.Range("D39").Value = '2.851.76'
$ws.Range("E39").Value = '  +4.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.795'
$ws.Range("E40").Value = '  -3.04%  '
$ws.Range("E41").Value = '  -4.25%  '
$ws.Range("E42").Value = '  -5.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0676'
$ws.Range("E43").Value = '  -2.48%  '
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.79'
$ws.Range("E45").Value = '  -4.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.39'
$ws.Range("E46").Value = '  -6.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '325.08'
$ws.Range("E47").Value = '  -5.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0273'
$ws.Range("E48").Value = '  -4.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '31.70'
$ws.Range("E49").Value = '  -4.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.992'
$ws.Range("E50").Value = '  -3.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.17'
$ws.Range("E51").Value = '  -2.91%  '
